$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.028.80"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.645.69"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "217.68"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "0.5174"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.06278"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "20.22"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D13").Value = "1.625.96"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "1.873.50"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "0.5547"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "0.0₅8083"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "64.90"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "26.012.54"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "4.582"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").Value = "190.36"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "5.869"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "144.18"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "0.1178"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").Value = "0.05319"
$ws.Range("E30").Value = "  -5.68%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "3.433"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "3.317"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "2.418"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "2.771"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "0.9369"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "0.5569"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "0.01567"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "5.754"
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "1.027.22"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").Value = "0.8220"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").Value = "100.52"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "1.783.21"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("D47").Value = "56.87"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "0.9977"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "0.4314"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "7.870"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "0.05112"
$ws.Range("E51").Value = "  -3.53%  "
